# Stand_Up_Meeting_Week_3.xlsx update
# "presentacion semana 4 y los stand up"
#
# Fills in the Week-3 stand-up answers for the first team member
# (rows 7-9, columns C-G) with the same per-member Q&A block pattern
# already used for the other team members further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - "Qué se hizo ayer?" answers
$ws.Range("C7").Value = "Nada"
$ws.Range("D7").Value = "Se cumplió lo propuesto "
$ws.Range("E7").Value = "Cumplí el objetivo"
$ws.Range("F7").Value = "Cumplí el objetivo"
$ws.Range("G7").Value = "Cumplí el objetivo"

# Row 8 - "Qué se hará hoy?" answers
$ws.Range("C8").Value = "Leer de las diapositivas del profe"
$ws.Range("D8").Value = "Asistir a la reunión para el diagrama de requisitos y ayudar con el diagrama de casos de uso"
$ws.Range("E8").Value = "Revisar las diapositivas"
$ws.Range("F8").Value = "Asistir a la presentación de avances en clase"
$ws.Range("G8").Value = "Asistir a reunión de asignación de tareas"

# Row 9 - "Qué cosas se oponen?" answers
$ws.Range("C9").Value = "Ninguna"
$ws.Range("D9").Value = "Nada"
$ws.Range("E9").Value = "Ninguna"
$ws.Range("F9").Value = "Nada"
$ws.Range("G9").Value = "Ninguna"

# F9 picks up the thin-bottom-border style used elsewhere in column F
# (e.g. F18) instead of the thicker divider style of its row neighbors.
$ws.Range("F9").Borders.Item(9).LineStyle = 1
$ws.Range("F9").Borders.Item(9).Weight = 2

# Update current selection to match the saved view state
$ws.Range("J13").Select() | Out-Null
